$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 848.63635
$ws.Range("I55").Value = 300.5
$ws.Range("J55").Value = 970.44446
$ws.Range("K55").Value = 300.5
$ws.Range("L55").Value = 970.44446
$ws.Range("M55").Value = -86.5
$ws.Range("N55").Value = -1398.44446
$ws.Range("H96").Value = 55562340
$ws.Range("I96").Value = 3526.3635
$ws.Range("J96").Value = 142869060
$ws.Range("K96").Value = 10579.0905
$ws.Range("L96").Value = 428607180
$ws.Range("M96").Value = -9206.0905
$ws.Range("N96").Value = -428609926
$ws.Range("H125").Value = 1189.2858
$ws.Range("I125").Value = 1028.4445
$ws.Range("J125").Value = 1478.8
$ws.Range("K125").Value = 9256.0005
$ws.Range("L125").Value = 13309.2
$ws.Range("M125").Value = -6796.0005
$ws.Range("N125").Value = -18229.2
$ws.Range("H129").Value = 18482.666
$ws.Range("I129").Value = 887.3158
$ws.Range("J129").Value = 27280.342
$ws.Range("K129").Value = 2661.9474
$ws.Range("L129").Value = 81841.026
$ws.Range("M129").Value = 2338.0526
$ws.Range("N129").Value = -91841.026
$ws.Range("H137").Value = 5508.0186
$ws.Range("I137").Value = 4608
$ws.Range("J137").Value = 7165.9473
$ws.Range("K137").Value = 13824
$ws.Range("L137").Value = 21497.8419
$ws.Range("M137").Value = -11274
$ws.Range("N137").Value = -26597.8419

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10350.386
$ws.Range("I32").Value = 9396.854499999999
$ws.Range("K32").Value = 9396.854499999999
$ws.Range("M32").Value = -9109.854499999999
$ws.Range("H61").Value = 2225.6
$ws.Range("I61").Value = 1372.8462
$ws.Range("J61").Value = 3809.2856
$ws.Range("K61").Value = 1372.8462
$ws.Range("L61").Value = 3809.2856
$ws.Range("M61").Value = -1160.8462
$ws.Range("N61").Value = -4233.2856
$ws.Range("H113").Value = 50000
$ws.Range("J113").Value = 50000
$ws.Range("L113").Value = 50000
$ws.Range("N113").Value = -58678
$ws.Range("H132").Value = 16669206
$ws.Range("I132").Value = 33336166
$ws.Range("J132").Value = 2245.5334
$ws.Range("K132").Value = 100008498
$ws.Range("L132").Value = 6736.600199999999
$ws.Range("M132").Value = -100005968
$ws.Range("N132").Value = -11796.6002
$ws.Range("H136").Value = 2225.6
$ws.Range("I136").Value = 1372.8462
$ws.Range("J136").Value = 3809.2856
$ws.Range("K136").Value = 4118.5386
$ws.Range("L136").Value = 11427.8568
$ws.Range("M136").Value = -1568.5386
$ws.Range("N136").Value = -16527.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 360268.1
$ws.Range("I134").Value = 1132.75
$ws.Range("J134").Value = 1274430.9
$ws.Range("K134").Value = 3398.25
$ws.Range("L134").Value = 3823292.7
$ws.Range("M134").Value = -863.25
$ws.Range("N134").Value = -3828362.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 35715320
$ws.Range("I2").Value = 83.333336
$ws.Range("J2").Value = 62501750
$ws.Range("K2").Value = 500.000016
$ws.Range("L2").Value = 375010500
$ws.Range("M2").Value = -387.000016
$ws.Range("N2").Value = -375010726
$ws.Range("H68").Value = 1322.4681
$ws.Range("I68").Value = 880.5
$ws.Range("J68").Value = 1649.8518
$ws.Range("K68").Value = 2641.5
$ws.Range("L68").Value = 4949.555399999999
$ws.Range("M68").Value = -1830.5
$ws.Range("N68").Value = -6571.555399999999
$ws.Range("H71").Value = 1322.4681
$ws.Range("I71").Value = 880.5
$ws.Range("J71").Value = 1649.8518
$ws.Range("K71").Value = 7924.5
$ws.Range("L71").Value = 14848.6662
$ws.Range("M71").Value = -3868.5
$ws.Range("N71").Value = -22960.6662
$ws.Range("H113").Value = 6421.3887
$ws.Range("I113").Value = 13078.75
$ws.Range("J113").Value = 1095.5
$ws.Range("K113").Value = 39236.25
$ws.Range("L113").Value = 3286.5
$ws.Range("M113").Value = -37066.25
$ws.Range("N113").Value = -7626.5
$ws.Range("H139").Value = 216823.22
$ws.Range("I139").Value = 233271.16
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 699813.48
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = -694673.48
$ws.Range("N139").Value = -19280
$ws.Range("H140").Value = 3107.4546
$ws.Range("I140").Value = 1600
$ws.Range("J140").Value = 4151.077
$ws.Range("K140").Value = 4800
$ws.Range("L140").Value = 12453.231
$ws.Range("M140").Value = 380
$ws.Range("N140").Value = -22813.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2044.1538
$ws.Range("I113").Value = 1973.4445
$ws.Range("J113").Value = 2203.25
$ws.Range("K113").Value = 1973.4445
$ws.Range("L113").Value = 2203.25
$ws.Range("M113").Value = 196.5554999999999
$ws.Range("N113").Value = -6543.25
$ws.Range("H127").Value = 41017
$ws.Range("J127").Value = 41017
$ws.Range("L127").Value = 41017
$ws.Range("N127").Value = -50937
$ws.Range("H130").Value = 49999.5
$ws.Range("J130").Value = 49999.5
$ws.Range("L130").Value = 49999.5
$ws.Range("N130").Value = -60039.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3931.7896
$ws.Range("J100").Value = 4800.4443
$ws.Range("L100").Value = 4800.4443
$ws.Range("N100").Value = -5882.4443
$ws.Range("H108").Value = 49626
$ws.Range("J108").Value = 49626
$ws.Range("L108").Value = 49626
$ws.Range("N108").Value = -57306
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("H130").Value = 44085.8
$ws.Range("J130").Value = 44085.8
$ws.Range("L130").Value = 44085.8
$ws.Range("N130").Value = -54125.8
$ws.Range("H132").Value = 4612.391
$ws.Range("I132").Value = 4392
$ws.Range("J132").Value = 4898.9
$ws.Range("K132").Value = 13176
$ws.Range("L132").Value = 14696.7
$ws.Range("M132").Value = -10646
$ws.Range("N132").Value = -19756.7
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H132").Value = 2119.2307
$ws.Range("I132").Value = 1810.45
$ws.Range("J132").Value = 2444.2632
$ws.Range("K132").Value = 5431.35
$ws.Range("L132").Value = 7332.7896
$ws.Range("M132").Value = -2901.35
$ws.Range("N132").Value = -12392.7896
$ws.Range("H138").Value = 45571.43
$ws.Range("J138").Value = 45571.43
$ws.Range("L138").Value = 45571.43
$ws.Range("N138").Value = -55851.43
$ws.Range("N125").ClearContents()
